$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Target cluster set to include "ECs" as a Sending cluster too (4x4 = 16 data rows)
# Re-write all data rows (2-17) with updated TPM-based values per commit "update scripts wuth new tpm"

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl12"
$ws.Range("C2").Value = "Ccr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 1.005755333333333
$ws.Range("H2").Value = 3.017266
$ws.Range("I2").Value = 0.01048729000197281
$ws.Range("J2").Value = 0.01048729000197281
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 1.454519666666667
$ws.Range("N2").Value = 4.363559
$ws.Range("O2").Value = 0.002411072527504041
$ws.Range("P2").Value = 0.002411072527504041
$ws.Range("Q2").Value = 1.462890912188223
$ws.Range("R2").Value = 13.166018209694
$ws.Range("S2").Value = 0.00002528561681172443
$ws.Range("T2").Value = 0.00002528561681172443

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl12"
$ws.Range("C3").Value = "Ccr1"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 1.005755333333333
$ws.Range("H3").Value = 3.017266
$ws.Range("I3").Value = 0.01048729000197281
$ws.Range("J3").Value = 0.01048729000197281
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 138.9276553333333
$ws.Range("N3").Value = 416.782966
$ws.Range("O3").Value = 0.2302922818860134
$ws.Range("P3").Value = 0.2302922818860133
$ws.Range("Q3").Value = 139.7272302989951
$ws.Range("R3").Value = 1257.545072690956
$ws.Range("S3").Value = 0.002415141945354691
$ws.Range("T3").Value = 0.002415141945354691

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ccl12"
$ws.Range("C4").Value = "Ccr1"
$ws.Range("D4").Value = "Neutrophils"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 1.005755333333333
$ws.Range("H4").Value = 3.017266
$ws.Range("I4").Value = 0.01048729000197281
$ws.Range("J4").Value = 0.01048729000197281
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 417.8699646666667
$ws.Range("N4").Value = 1253.609894
$ws.Range("O4").Value = 0.6926787000314772
$ws.Range("P4").Value = 0.6926787000314772
$ws.Range("Q4").Value = 420.2749456033116
$ws.Range("R4").Value = 3782.474510429804
$ws.Range("S4").Value = 0.007264322405419631
$ws.Range("T4").Value = 0.007264322405419631

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Ccl12"
$ws.Range("C5").Value = "Ccr1"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 1.005755333333333
$ws.Range("H5").Value = 3.017266
$ws.Range("I5").Value = 0.01048729000197281
$ws.Range("J5").Value = 0.01048729000197281
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 45.01451866666667
$ws.Range("N5").Value = 135.043556
$ws.Range("O5").Value = 0.07461794555500532
$ws.Range("P5").Value = 0.07461794555500531
$ws.Range("Q5").Value = 45.27359222643289
$ws.Range("R5").Value = 407.462330037896
$ws.Range("S5").Value = 0.0007825400343867585
$ws.Range("T5").Value = 0.0007825400343867583

# Row 6
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Ccl12"
$ws.Range("C6").Value = "Ccr1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 38.755371
$ws.Range("H6").Value = 116.266113
$ws.Range("I6").Value = 0.4041130097356814
$ws.Range("J6").Value = 0.4041130097356814
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 1.454519666666667
$ws.Range("N6").Value = 4.363559
$ws.Range("O6").Value = 0.002411072527504041
$ws.Range("P6").Value = 0.002411072527504041
$ws.Range("Q6").Value = 56.370449308463
$ws.Range("R6").Value = 507.334043776167
$ws.Range("S6").Value = 0.0009743457757806745
$ws.Range("T6").Value = 0.0009743457757806742

# Row 7
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Ccl12"
$ws.Range("C7").Value = "Ccr1"
$ws.Range("D7").Value = "Inflammatory-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 38.755371
$ws.Range("H7").Value = 116.266113
$ws.Range("I7").Value = 0.4041130097356814
$ws.Range("J7").Value = 0.4041130097356814
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 138.9276553333333
$ws.Range("N7").Value = 416.782966
$ws.Range("O7").Value = 0.2302922818860134
$ws.Range("P7").Value = 0.2302922818860133
$ws.Range("Q7").Value = 5384.192824603461
$ws.Range("R7").Value = 48457.73542143115
$ws.Range("S7").Value = 0.0930641071518548
$ws.Range("T7").Value = 0.09306410715185479

# Row 8
$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("B8").Value = "Ccl12"
$ws.Range("C8").Value = "Ccr1"
$ws.Range("D8").Value = "Neutrophils"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 38.755371
$ws.Range("H8").Value = 116.266113
$ws.Range("I8").Value = 0.4041130097356814
$ws.Range("J8").Value = 0.4041130097356814
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 417.8699646666667
$ws.Range("N8").Value = 1253.609894
$ws.Range("O8").Value = 0.6926787000314772
$ws.Range("P8").Value = 0.6926787000314772
$ws.Range("Q8").Value = 16194.70551041356
$ws.Range("R8").Value = 145752.349593722
$ws.Range("S8").Value = 0.2799204742495195
$ws.Range("T8").Value = 0.2799204742495195

# Row 9
$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("B9").Value = "Ccl12"
$ws.Range("C9").Value = "Ccr1"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 38.755371
$ws.Range("H9").Value = 116.266113
$ws.Range("I9").Value = 0.4041130097356814
$ws.Range("J9").Value = 0.4041130097356814
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 45.01451866666667
$ws.Range("N9").Value = 135.043556
$ws.Range("O9").Value = 0.07461794555500532
$ws.Range("P9").Value = 0.07461794555500531
$ws.Range("Q9").Value = 1744.554371313092
$ws.Range("R9").Value = 15700.98934181783
$ws.Range("S9").Value = 0.03015408255852641
$ws.Range("T9").Value = 0.0301540825585264

# Row 10
$ws.Range("A10").Value = "Neutrophils"
$ws.Range("B10").Value = "Ccl12"
$ws.Range("C10").Value = "Ccr1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.988471999999999
$ws.Range("H10").Value = 23.965416
$ws.Range("I10").Value = 0.08329801469605898
$ws.Range("J10").Value = 0.08329801469605898
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 1.454519666666667
$ws.Range("N10").Value = 4.363559
$ws.Range("O10").Value = 0.002411072527504041
$ws.Range("P10").Value = 0.002411072527504041
$ws.Range("Q10").Value = 11.619389630616
$ws.Range("R10").Value = 104.574506675544
$ws.Range("S10").Value = 0.0002008375548292957
$ws.Range("T10").Value = 0.0002008375548292956

# Row 11
$ws.Range("A11").Value = "Neutrophils"
$ws.Range("B11").Value = "Ccl12"
$ws.Range("C11").Value = "Ccr1"
$ws.Range("D11").Value = "Inflammatory-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 7.988471999999999
$ws.Range("H11").Value = 23.965416
$ws.Range("I11").Value = 0.08329801469605898
$ws.Range("J11").Value = 0.08329801469605898
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 138.9276553333333
$ws.Range("N11").Value = 416.782966
$ws.Range("O11").Value = 0.2302922818860134
$ws.Range("P11").Value = 0.2302922818860133
$ws.Range("Q11").Value = 1109.819684655984
$ws.Range("R11").Value = 9988.377161903854
$ws.Range("S11").Value = 0.0191828898809301
$ws.Range("T11").Value = 0.01918288988093009

# Row 12
$ws.Range("A12").Value = "Neutrophils"
$ws.Range("B12").Value = "Ccl12"
$ws.Range("C12").Value = "Ccr1"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 7.988471999999999
$ws.Range("H12").Value = 23.965416
$ws.Range("I12").Value = 0.08329801469605898
$ws.Range("J12").Value = 0.08329801469605898
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 417.8699646666667
$ws.Range("N12").Value = 1253.609894
$ws.Range("O12").Value = 0.6926787000314772
$ws.Range("P12").Value = 0.6926787000314772
$ws.Range("Q12").Value = 3338.142512380656
$ws.Range("R12").Value = 30043.2826114259
$ws.Range("S12").Value = 0.05769876053486902
$ws.Range("T12").Value = 0.05769876053486902

# Row 13
$ws.Range("A13").Value = "Neutrophils"
$ws.Range("B13").Value = "Ccl12"
$ws.Range("C13").Value = "Ccr1"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 7.988471999999999
$ws.Range("H13").Value = 23.965416
$ws.Range("I13").Value = 0.08329801469605898
$ws.Range("J13").Value = 0.08329801469605898
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 45.01451866666667
$ws.Range("N13").Value = 135.043556
$ws.Range("O13").Value = 0.07461794555500532
$ws.Range("P13").Value = 0.07461794555500531
$ws.Range("Q13").Value = 359.597221962144
$ws.Range("R13").Value = 3236.374997659296
$ws.Range("S13").Value = 0.006215526725430562
$ws.Range("T13").Value = 0.006215526725430561

# Row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Ccl12"
$ws.Range("C14").Value = "Ccr1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 48.15271133333334
$ws.Range("H14").Value = 144.458134
$ws.Range("I14").Value = 0.5021016855662869
$ws.Range("J14").Value = 0.5021016855662868
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 1.454519666666667
$ws.Range("N14").Value = 4.363559
$ws.Range("O14").Value = 0.002411072527504041
$ws.Range("P14").Value = 0.002411072527504041
$ws.Range("Q14").Value = 70.03906563765624
$ws.Range("R14").Value = 630.3515907389061
$ws.Range("S14").Value = 0.001210603580082347
$ws.Range("T14").Value = 0.001210603580082346

# Row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Ccl12"
$ws.Range("C15").Value = "Ccr1"
$ws.Range("D15").Value = "Inflammatory-Mac"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 48.15271133333334
$ws.Range("H15").Value = 144.458134
$ws.Range("I15").Value = 0.5021016855662869
$ws.Range("J15").Value = 0.5021016855662868
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 138.9276553333333
$ws.Range("N15").Value = 416.782966
$ws.Range("O15").Value = 0.2302922818860134
$ws.Range("P15").Value = 0.2302922818860133
$ws.Range("Q15").Value = 6689.743283482827
$ws.Range("R15").Value = 60207.68955134544
$ws.Range("S15").Value = 0.1156301429078738
$ws.Range("T15").Value = 0.1156301429078737

# Row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Ccl12"
$ws.Range("C16").Value = "Ccr1"
$ws.Range("D16").Value = "Neutrophils"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 48.15271133333334
$ws.Range("H16").Value = 144.458134
$ws.Range("I16").Value = 0.5021016855662869
$ws.Range("J16").Value = 0.5021016855662868
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 417.8699646666667
$ws.Range("N16").Value = 1253.609894
$ws.Range("O16").Value = 0.6926787000314772
$ws.Range("P16").Value = 0.6926787000314772
$ws.Range("Q16").Value = 20121.5717834642
$ws.Range("R16").Value = 181094.1460511778
$ws.Range("S16").Value = 0.3477951428416691
$ws.Range("T16").Value = 0.3477951428416691

# Row 17
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Ccl12"
$ws.Range("C17").Value = "Ccr1"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 48.15271133333334
$ws.Range("H17").Value = 144.458134
$ws.Range("I17").Value = 0.5021016855662869
$ws.Range("J17").Value = 0.5021016855662868
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 45.01451866666667
$ws.Range("N17").Value = 135.043556
$ws.Range("O17").Value = 0.07461794555500532
$ws.Range("P17").Value = 0.07461794555500531
$ws.Range("Q17").Value = 2167.571123164945
$ws.Range("R17").Value = 19508.1401084845
$ws.Range("S17").Value = 0.0374657962366616
$ws.Range("T17").Value = 0.03746579623666158
